# NYPD CompStat weekly report update — new crime data collected for the
# report covering 12/9/2024 through 12/15/2024 (Volume 31, Number 50).
#
# 1) Update the report header: issue number and reporting week dates.
# 2) Widen column H slightly to fit the new "28 Day % Chg" values.
# 3) Refresh the crime-complaint statistics table (rows 14-31): several
#    categories that previously had no complaints ("0"/"***.*" placeholder
#    text) now have real counts, so their cells are converted from text to
#    numeric/percentage, and all of the counts/percentages are refreshed
#    with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  49" -> "...Number  50" ---------------
$ws.Range("A8").Characters(21,2).Text = "50"

# --- Header: reporting week dates ---------------------------------------
$ws.Range("C9").Characters(27,9).Text = "12/9/2024"
$ws.Range("C9").Characters(47,9).Text = "12/15/2024"

# --- Column H is now a bit wider to fit the refreshed values -----------
$ws.Columns.Item(8).ColumnWidth = 7.433768

# --- Cells that go from the "no data" placeholder text to real numeric/
#     percentage data need their number format set explicitly so they
#     match the sheet's existing numeric ("#,##0") and percentage
#     ("#,##0.0;"-"#,##0.0") styles used elsewhere in the table. ---------
    $ws.Range("D15").NumberFormat = '#,##0'
    $ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G15").NumberFormat = '#,##0'
    $ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("C18").NumberFormat = '#,##0'
    $ws.Range("D20").NumberFormat = '#,##0'
    $ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("D27").NumberFormat = '#,##0'
    $ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G27").NumberFormat = '#,##0'
    $ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("C28").NumberFormat = '#,##0'
    $ws.Range("D31").NumberFormat = '#,##0'
    $ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G31").NumberFormat = '#,##0'
    $ws.Range("H31").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Refreshed crime-complaint figures (rows 14-31) ---------------------
    # Row 14
    $ws.Range("M14").Value = -38.461538461538
    $ws.Range("N14").Value = -77.142857142857

    # Row 15
    $ws.Range("D15").Value = 1
    $ws.Range("E15").Value = -100
    $ws.Range("F15").Value = 1
    $ws.Range("G15").Value = 1
    $ws.Range("H15").Value = 0
    $ws.Range("J15").Value = 12
    $ws.Range("K15").Value = 33.333333333333

    # Row 16
    $ws.Range("C16").Value = 3
    $ws.Range("D16").Value = 5
    $ws.Range("E16").Value = -40
    $ws.Range("F16").Value = 16
    $ws.Range("G16").Value = 17
    $ws.Range("H16").Value = -5.882352941176
    $ws.Range("I16").Value = 249
    $ws.Range("J16").Value = 269
    $ws.Range("K16").Value = -7.434944237918
    $ws.Range("L16").Value = 9.210526315789
    $ws.Range("M16").Value = 12.162162162162
    $ws.Range("N16").Value = -71.800679501698

    # Row 17
    $ws.Range("C17").Value = 7
    $ws.Range("D17").Value = 7
    $ws.Range("E17").Value = 0
    $ws.Range("F17").Value = 24
    $ws.Range("G17").Value = 28
    $ws.Range("H17").Value = -14.285714285714
    $ws.Range("I17").Value = 464
    $ws.Range("J17").Value = 374
    $ws.Range("K17").Value = 24.064171122994
    $ws.Range("L17").Value = 7.407407407407
    $ws.Range("M17").Value = 109.009009009009
    $ws.Range("N17").Value = -28.834355828220

    # Row 18
    $ws.Range("C18").Value = 1
    $ws.Range("D18").Value = 3
    $ws.Range("E18").Value = -66.666666666666
    $ws.Range("G18").Value = 9
    $ws.Range("H18").Value = -55.555555555555
    $ws.Range("I18").Value = 98
    $ws.Range("J18").Value = 115
    $ws.Range("K18").Value = -14.782608695652
    $ws.Range("L18").Value = -45.555555555555
    $ws.Range("M18").Value = -5.769230769230
    $ws.Range("N18").Value = -80.970873786407

    # Row 19
    $ws.Range("C19").Value = 7
    $ws.Range("D19").Value = 8
    $ws.Range("E19").Value = -12.5
    $ws.Range("F19").Value = 28
    $ws.Range("G19").Value = 32
    $ws.Range("H19").Value = -12.5
    $ws.Range("I19").Value = 472
    $ws.Range("J19").Value = 534
    $ws.Range("K19").Value = -11.610486891385
    $ws.Range("L19").Value = -21.464226289517
    $ws.Range("M19").Value = 108.849557522124
    $ws.Range("N19").Value = 24.538258575197

    # Row 20
    $ws.Range("C20").Value = 1
    $ws.Range("D20").Value = 1
    $ws.Range("E20").Value = 0
    $ws.Range("F20").Value = 8
    $ws.Range("G20").Value = 3
    $ws.Range("H20").Value = 166.666666666667
    $ws.Range("I20").Value = 107
    $ws.Range("J20").Value = 84
    $ws.Range("K20").Value = 27.380952380952
    $ws.Range("L20").Value = 25.882352941176
    $ws.Range("M20").Value = 57.352941176470
    $ws.Range("N20").Value = -67.278287461773

    # Row 21
    $ws.Range("D21").Value = 25
    $ws.Range("E21").Value = -24
    $ws.Range("F21").Value = 82
    $ws.Range("G21").Value = 90
    $ws.Range("H21").Value = -8.888888888888
    $ws.Range("I21").Value = 1414
    $ws.Range("J21").Value = 1395
    $ws.Range("K21").Value = 1.362007168458
    $ws.Range("L21").Value = -8.538163001293
    $ws.Range("M21").Value = 62.155963302752
    $ws.Range("N21").Value = -50.052984811020

    # Row 22
    $ws.Range("G22").Value = 1
    $ws.Range("H22").Value = 0
    $ws.Range("L22").Value = -60.975609756097

    # Row 23
    $ws.Range("C23").Value = 2
    $ws.Range("D23").Value = 4
    $ws.Range("E23").Value = -50
    $ws.Range("F23").Value = 12
    $ws.Range("H23").Value = 20
    $ws.Range("I23").Value = 148
    $ws.Range("J23").Value = 157
    $ws.Range("K23").Value = -5.732484076433
    $ws.Range("L23").Value = -17.318435754189
    $ws.Range("M23").Value = 48

    # Row 24
    $ws.Range("C24").Value = 20
    $ws.Range("D24").Value = 21
    $ws.Range("E24").Value = -4.761904761904
    $ws.Range("F24").Value = 62
    $ws.Range("G24").Value = 64
    $ws.Range("H24").Value = -3.125
    $ws.Range("I24").Value = 930
    $ws.Range("J24").Value = 1244
    $ws.Range("K24").Value = -25.241157556270
    $ws.Range("L24").Value = -29.438543247344
    $ws.Range("M24").Value = -1.587301587301

    # Row 25
    $ws.Range("C25").Value = 5
    $ws.Range("D25").Value = 8
    $ws.Range("E25").Value = -37.5
    $ws.Range("G25").Value = 19
    $ws.Range("H25").Value = -52.631578947368
    $ws.Range("I25").Value = 181
    $ws.Range("J25").Value = 550
    $ws.Range("K25").Value = -67.090909090909
    $ws.Range("L25").Value = -75.769745649263

    # Row 26
    $ws.Range("C26").Value = 7
    $ws.Range("D26").Value = 10
    $ws.Range("E26").Value = -30
    $ws.Range("F26").Value = 27
    $ws.Range("G26").Value = 48
    $ws.Range("H26").Value = -43.75
    $ws.Range("I26").Value = 704
    $ws.Range("J26").Value = 528
    $ws.Range("K26").Value = 33.333333333333
    $ws.Range("L26").Value = 31.098696461825
    $ws.Range("M26").Value = 47.280334728033

    # Row 27
    $ws.Range("D27").Value = 1
    $ws.Range("E27").Value = -100
    $ws.Range("F27").Value = 1
    $ws.Range("G27").Value = 1
    $ws.Range("H27").Value = 0
    $ws.Range("J27").Value = 22
    $ws.Range("K27").Value = -4.545454545454

    # Row 28
    $ws.Range("C28").Value = 1
    $ws.Range("D28").Value = 5
    $ws.Range("E28").Value = -80
    $ws.Range("F28").Value = 4
    $ws.Range("G28").Value = 10
    $ws.Range("H28").Value = -60
    $ws.Range("I28").Value = 60
    $ws.Range("J28").Value = 72
    $ws.Range("K28").Value = -16.666666666666
    $ws.Range("L28").Value = -3.225806451612

    # Row 29
    $ws.Range("F29").Value = 1
    $ws.Range("M29").Value = -60
    $ws.Range("N29").Value = -78.947368421052

    # Row 30
    $ws.Range("F30").Value = 1
    $ws.Range("M30").Value = -59.375
    $ws.Range("N30").Value = -81.690140845070

    # Row 31
    $ws.Range("D31").Value = 1
    $ws.Range("E31").Value = -100
    $ws.Range("G31").Value = 1
    $ws.Range("H31").Value = -100
    $ws.Range("J31").Value = 5
    $ws.Range("K31").Value = 20

